# Update "想去人数" (want-to-go count) values that changed in the latest
# gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 378
$wsExhibit.Range("F3").Value = 831
$wsExhibit.Range("F5").Value = 1025
$wsExhibit.Range("F6").Value = 2392
$wsExhibit.Range("F7").Value = 202

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 378
$wsAll.Range("F3").Value = 831
$wsAll.Range("F7").Value = 1025
$wsAll.Range("F8").Value = 2392
$wsAll.Range("F10").Value = 202
